{"js": "// The edit deletes four paragraphs that directly follow the paragraph\n// containing \"LOB1036: Geometria Anal\u00edtica (Requisito fraco)\":\n//   1) an empty paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) another empty paragraph\n//   4) an empty paragraph that carries <w:pageBreakBefore/> and jc=\"left\"\n//\n// Find the anchor paragraph by its text, then remove the next four\n// paragraphs (in reverse order, so earlier deletes don't shift later ones).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst anchorText = \"LOB1036: Geometria Anal\u00edtica (Requisito fraco)\";\n\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// Delete the four paragraphs immediately after the anchor paragraph.\n// Delete from the highest index down to keep the lower indices valid.\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i <= anchorIndex + 4 && i < items.length; i++) {\n  toDelete.push(items[i]);\n}\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# The edit removes four paragraphs that directly follow the paragraph\n# containing \"LOB1036: Geometria Anal\u00edtica (Requisito fraco)\":\n#   1) an empty paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) another empty paragraph\n#   4) an empty paragraph that carries a page-break-before and jc=\"left\"\n#\n# Locate the anchor paragraph by its text, then delete the range that spans\n# the four paragraphs immediately following it (start of the first one to\n# the end of the fourth one, including their paragraph marks).\n\n$d = $word.ActiveDocument\n\n$anchorText = \"LOB1036: Geometria Anal\u00edtica (Requisito fraco)\"\n\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    $t = $t.Trim()\n    if ($t -eq $anchorText) {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -eq $null) {\n    throw \"Could not find anchor paragraph: $anchorText\"\n}\n\n$firstToRemove = $anchor.Next()\n$lastToRemove = $firstToRemove\nfor ($i = 1; $i -lt 4; $i++) {\n    $lastToRemove = $lastToRemove.Next()\n}\n\n$rng = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)\n$rng.Delete()\n"}
